$d = $word.ActiveDocument

$pairs = @(
    @("535÷4=133, 3", "834÷5=166, 4"),
    @("697÷5=139, 2", "613÷4=153, 1"),
    @("330÷4=82, 2", "930÷3=310, 0"),
    @("841÷8=105, 1", "555÷4=138, 3"),
    @("901÷9=100, 1", "349÷3=116, 1"),
    @("356÷6=59, 2", "965÷3=321, 2"),
    @("209÷3=69, 2", "395÷5=79, 0"),
    @("869÷7=124, 1", "135÷4=33, 3"),
    @("212÷7=30, 2", "753÷3=251, 0"),
    @("791÷5=158, 1", "216÷9=24, 0"),
    @("761÷9=84, 5", "704÷7=100, 4"),
    @("738÷5=147, 3", "959÷2=479, 1"),
    @("833÷6=138, 5", "157÷8=19, 5"),
    @("956÷8=119, 4", "120÷4=30, 0"),
    @("540÷9=60, 0", "276÷9=30, 6"),
    @("102÷2=51, 0", "646÷3=215, 1"),
    @("701÷6=116, 5", "858÷3=286, 0"),
    @("501÷5=100, 1", "754÷3=251, 1"),
    @("721÷4=180, 1", "269÷6=44, 5"),
    @("399÷4=99, 3", "699÷7=99, 6"),
    @("810÷3=270, 0", "969÷3=323, 0"),
    @("441÷3=147, 0", "284÷3=94, 2"),
    @("831÷3=277, 0", "366÷3=122, 0"),
    @("719÷8=89, 7", "626÷4=156, 2"),
    @("849÷3=283, 0", "707÷7=101, 0")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
